$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold an exact text value, regardless of whether the
# string would otherwise be auto-converted by Excel into a number/date,
# then restore the default (unstyled) cell style so no stray formatting
# is introduced.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.190.73"
$ws.Range("E2").Value = "  +0.31%  "
Set-TextValue $ws.Range("D3") "1.887.49"
$ws.Range("E3").Value = "  -0.55%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "322.46"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.08006"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "0.9923"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D11") "22.61"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D12") "5.916"
$ws.Range("E12").Value = "  -0.38%  "
Set-TextValue $ws.Range("D13") "1.843.71"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "7.015"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "89.32"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.002"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.06630"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001024"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "17.38"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "1.001"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D21") "29.194.97"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.487"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D23") "11.67"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.176"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D25") "2.118.97"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "154.84"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "19.60"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "6.002"
$ws.Range("E28").Value = "  +6.93%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.079"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "117.00"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "1.023"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.09401"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D33") "3.539"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "1.376"
$ws.Range("E34").Value = "  -3.04%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D35") "5.348"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.06038"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02234"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.170"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D39") "7.985"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.5806"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D41") "0.1826"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "10.01"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "2.426"
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.273"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D45") "0.07700"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "12.13"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5455"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.899"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "113.22"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue $ws.Range("D50") "0.2939"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D51") "43.59"
$ws.Range("E51").Value = "  -1.78%  "
